# Update the build timestamp embedded in the version strings throughout the
# workbook: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" worksheet ---
$about = $wb.Worksheets.Item("About")

$a2 = $about.Range("A2")
$a2text = $a2.Value()
$a2.Value = $a2text.Replace($oldStamp, $newStamp)

$a6 = $about.Range("A6")
$a6text = $a6.Value()
$a6.Value = $a6text.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" worksheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 9; $row++) {
    $cell = $data.Cells.Item($row, 19)  # column S
    $cellText = $cell.Value()
    $cell.Value = $cellText.Replace($oldStamp, $newStamp)
}
